$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting D:K -> E:L
$ws.Columns("D").Insert()

# Propagate formatting (number format / font / style) from the former D column
# (now shifted to E) into the newly inserted, blank D column so every row's
# new D cell matches its row-type style (date header rows vs numeric rows).
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 5, 6, 37 and 79 are section-label rows that only have content in
# column A/B and never had a D:K value - remove the blank D cell that the
# format-paste introduced so they stay cell-less in column D, same as before.
$ws.Range("D5").Clear()
$ws.Range("D6").Clear()
$ws.Range("D37").Clear()
$ws.Range("D79").Clear()

# Fill in the new column D with the latest reporting period's figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 15668200
$ws.Range("D9").Value = 10775300
$ws.Range("D10").Value = 4892900
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 93400
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 13750800
$ws.Range("D18").Value = 1917400
$ws.Range("D20").Value = -36400
$ws.Range("D21").Value = 2242500
$ws.Range("D22").Value = 220700
$ws.Range("D23").Value = 1660300
$ws.Range("D24").Value = 290300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 1370000
$ws.Range("D27").Value = 1350100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -12500
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 36400
$ws.Range("D33").Value = 1337600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 1337600
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 903400
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 2679200
$ws.Range("D44").Value = 1677800
$ws.Range("D45").Value = 471600
$ws.Range("D46").Value = 5732000
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 1730800
$ws.Range("D49").Value = 9594200
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 857900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 17914900
$ws.Range("D57").Value = 1705300
$ws.Range("D58").Value = 350600
$ws.Range("D59").Value = 2259800
$ws.Range("D60").Value = 4315700
$ws.Range("D61").Value = 3740700
$ws.Range("D62").Value = 2793700
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 10892200
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 9439800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 7022700
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 1337600
$ws.Range("D83").Value = 361500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 1407800
$ws.Range("D91").Value = -365600
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -629400
$ws.Range("D96").Value = -479500
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -1378800
$ws.Range("D101").Value = -45600
$ws.Range("D102").Value = -646000
